# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# sheets, which carry the same data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 11732
    "F3"  = 11362
    "F6"  = 1028
    "F11" = 10774
    "F16" = 2469
    "F22" = 10928
    "F24" = 32
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
